$d = $word.ActiveDocument

# Locate the paragraph holding the "validity period" sentence. It currently
# reads (field results shown in guillemets):
#   ... tj. z dniem «validity_date» r. i obowiązuje do dnia «validity_date_end» r.
# The edit collapses " i obowiązuje do dnia " to a single space and drops the
# now-redundant trailing " r." after the closing mergefield, leaving:
#   ... tj. z dniem «validity_date» r. «validity_date_end»
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains("obowiązuje do dnia")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # 1) Drop the trailing " r." run (the very last 3 characters of the
    #    paragraph, right before the paragraph mark) that used to follow the
    #    validity_date_end mergefield / _GoBack bookmark.
    $paraEnd = $target.Range.End
    $tail = $d.Range($paraEnd - 4, $paraEnd - 1)
    if ($tail.Text -eq " r.") {
        $tail.Delete()
    }
}

# 2) " i obowiązuje do dnia " -> " " (collapse the run text to a single space)
$d.Content.Find.Execute(
    " i obowiązuje do dnia ", $true, $false, $false, $false, $false,
    $true, 1, $false, " ", 2
) | Out-Null
